$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "DONE" / "?" status markers in column E first (matches
# shared-string insertion order: DONE then ? are the first two new strings)
$ws.Range("E3").Value = "DONE"
$ws.Range("E5").Value = "DONE"
$ws.Range("E7").Value = "x"
$ws.Range("E8").Value = "DONE"
$ws.Range("E9").Value = "?"
$ws.Range("E11").Value = "?"
$ws.Range("E12").Value = "DONE"
$ws.Range("E14").Value = "x"

# New "First Milestone" column header and notes
$ws.Range("F1").Value = "First Milestone"
$ws.Range("F3").Value = "x"
$ws.Range("F4").Value = "x"
$ws.Range("F5").Value = "player controls"
$ws.Range("F6").Value = "obstacles generating"
$ws.Range("F2").Value = "60 Points Scaled"

# Set column F to a best-fit width sized for its new content (e.g. "obstacles generating")
$ws.Range("F1").EntireColumn.ColumnWidth = 11.4

# Update the active cell selection
$ws.Range("H8").Select()
